$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.161.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.91%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.14%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.538"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.24%  "

# Row 9
$ws.Range("E9").Value = "  +6.73%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.05"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0819"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.64%  "

# Row 12
$ws.Range("E12").Value = "  +0.90%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.663.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.300.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.94%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.810"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.066.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.92%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.49%  "

# Row 20
$ws.Range("E20").Value = "  +2.51%  "

# Row 21
$ws.Range("E21").Value = "  +2.39%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.76%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.48%  "

# Row 24
$ws.Range("E24").Value = "  +2.63%  "

# Row 25
$ws.Range("E25").Value = "  +3.88%  "

# Row 26
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.82%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.22%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.74%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.82"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("E34").Value = "  -1.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0742"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.71%  "

# Row 38
$ws.Range("E38").Value = "  +0.94%  "

# Row 39
$ws.Range("E39").Value = "  +1.83%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.07%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.72%  "

# Row 42
$ws.Range("E42").Value = "  +0.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.982.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.08"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.78%  "

# Row 46
$ws.Range("E46").Value = "  +3.73%  "

# Row 47
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +19.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.531.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.76%  "

# Row 51
$ws.Range("E51").Value = "  +2.23%  "
